$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each cell update from the source diff. Column D price cells are
# forced to Text (NumberFormat "@") before the write so plain-numeric-looking
# strings (e.g. "543.17") are not auto-converted to numbers by Excel, matching
# the original inlineStr cell type. The style is then reset to "Normal" so no
# extra/visible formatting is left behind on the cell.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.894.21'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.54%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.354.38'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.50%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '543.17'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.50'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.49%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.570'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +6.11%  '
$ws.Range('E9').Value = '  -0.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.52'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.62%  '
$ws.Range('E11').Value = '  -1.51%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.357'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.773.47'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.57%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.72'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '57.829.43'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.57%  '
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.354.59'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.94%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.73'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.66%  '
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.29'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.40%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '329.94'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.61%  '
$ws.Range('E21').Value = '  -2.28%  '
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '62.51'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.49%  '
$ws.Range('E24').Value = '  -1.81%  '
$ws.Range('E25').Value = '  +0.22%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.40'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.58%  '
$ws.Range('E27').Value = '  -3.45%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '170.34'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.62%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0736'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.13'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.91%  '
$ws.Range('E32').Value = '  +0.18%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.38'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.91%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.22'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.44%  '
$ws.Range('E36').Value = '  +0.22%  '
$ws.Range('E37').Value = '  -2.12%  '
$ws.Range('E38').Value = '  -0.68%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '39.06'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '142.61'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.69%  '
$ws.Range('E41').Value = '  -0.20%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.65'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.36%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '289.25'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.05%  '
$ws.Range('E44').Value = '  +2.02%  '
$ws.Range('E45').Value = '  +0.64%  '
$ws.Range('E46').Value = '  -1.30%  '
$ws.Range('E47').Value = '  +0.88%  '
$ws.Range('E48').Value = '  +1.68%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.382'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.16%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.49'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.08'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.51%  '
